$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded between the existing rows
# 420 and 421 (by date order in the source feed). Insert a fresh row at
# position 421, which pushes the former rows 421-516 down to 422-517 and
# grows the used range to A1:R517, then populate the new row with its data.
$ws.Rows(421).Insert()

$ws.Range("A421").Value = 3
$ws.Range("B421").Value = "Femacal de La Calera"
$ws.Range("C421").Value = "Coquimbo"
$ws.Range("D421").Value = 45173
$ws.Range("E421").Value = 5
$ws.Range("F421").Value = 100112001
$ws.Range("G421").Value = "Berenjena"
$ws.Range("H421").Value = "Sin especificar"
$ws.Range("I421").Value = "Primera"
$ws.Range("J421").Value = 40
$ws.Range("K421").Value = 9000
$ws.Range("L421").Value = 9000
$ws.Range("M421").Value = 9000
$ws.Range("N421").Value = "$/caja 60 unidades"
$ws.Range("O421").Value = "Región de Arica y Parinacota"
$ws.Range("P421").Value = 150
$ws.Range("Q421").Value = 60
$ws.Range("R421").Value = "Hortaliza"
